# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet (5th sheet), filling every existing data row
# with the report date, the legislator's name and id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$lastRow = 33
$dateValue = "2011-11-22"
$legislatorName = "徐少萍"
$legislatorId = 726

# ---- Header row (row 1): copy formatting from an existing header cell so
# the new header cells (H1:J1) get the same bold/border style used by the
# rest of the header row, then set their labels.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# ---- Data rows (2..33): the date column must stay plain text (otherwise
# Excel silently turns the ISO-looking string into a date serial number),
# so mark the column as Text before writing the values.
$ws.Range("H2:H" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $dateValue
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
